$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.417.66'
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.563.55'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.64%  '
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '285.31'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3637'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.75%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.57'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.63%  '
$ws.Range("E9").Value = '  -1.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.129'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07413'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.91%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.83'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.933'
$ws.Range("D14").ClearFormats()
$ws.Range("E15").Value = '  -0.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.563.75'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001105'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '88.20'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06686'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.74%  '
$ws.Range("E20").Value = '  -0.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.365'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.13'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.00'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.408.81'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.409'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.559'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '149.56'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.64%  '
$ws.Range("E28").Value = '  -3.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.998'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.13'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.739.15'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.064'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.148'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.994'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.64%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.832'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("E36").Value = '  -1.84%  '
$ws.Range("E37").Value = '  -2.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.306'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06388'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2210'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.338'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6091'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9998'
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.83'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.756'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5764'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.94%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.013'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.72'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.216'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07212'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.52%  '
